# Daily injury/props sheet refresh: the block of player rows 6-11 rotates by
# one slot. D'Angelo Russell's row (previously last, at row 11) moves up to
# the top of the block at row 6, and every other row shifts down by one
# (old row 6 -> new row 7, old row 7 -> new row 8, ... old row 10 -> new
# row 11). Russell's injury status also changes from "Out" to "Day-To-Day".
#
# Columns E, U and V are blank for every row in this block both before and
# after the edit, so they are left untouched. Column C ("Statut") is blank
# for every row except row 11 ("Out") before the edit, and blank for every
# row except the new row 6 ("Day-To-Day") after - so only C6 and C11 need an
# explicit write there; the rest of column C is already correct and is left
# alone.
#
# Capture each row's data (split around the untouched blank columns) before
# writing anything back, so the cyclic rotation doesn't clobber a row we
# still need to read from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowSegments($row) {
    return @(
        $ws.Range("A$($row):B$($row)").Value2,
        $ws.Range("D$($row)").Value2,
        $ws.Range("F$($row):T$($row)").Value2,
        $ws.Range("W$($row):AL$($row)").Value2
    )
}

function Set-RowSegments($row, $segs) {
    $ws.Range("A$($row):B$($row)").Value = $segs[0]
    $ws.Range("D$($row)").Value = $segs[1]
    $ws.Range("F$($row):T$($row)").Value = $segs[2]
    $ws.Range("W$($row):AL$($row)").Value = $segs[3]
}

$seg6  = Get-RowSegments 6
$seg7  = Get-RowSegments 7
$seg8  = Get-RowSegments 8
$seg9  = Get-RowSegments 9
$seg10 = Get-RowSegments 10
$seg11 = Get-RowSegments 11

# Old row 11 (D'Angelo Russell) becomes the new row 6; everyone else shifts
# down by one row.
Set-RowSegments 6  $seg11
Set-RowSegments 7  $seg6
Set-RowSegments 8  $seg7
Set-RowSegments 9  $seg8
Set-RowSegments 10 $seg9
Set-RowSegments 11 $seg10

# Status column: Russell is now "Day-To-Day" at the top of the block, and
# his old "Out" at row 11 is gone (blank, like the rest of the column).
$ws.Range("C6").Value = "Day-To-Day"
$ws.Range("C11").Value = ""
